$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @(3.230985683306322, 1.667794583268128, 26.21740644021617, 8.660232485948974, 39.7764191927396)
    3  = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671)
    4  = @(0.127881588408715, 0.3127903958511391, 3.900430680208489, 0.496779210170732, 4.837881874639075)
    5  = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    6  = @(1.459612070389937, 1.667794583268128, 3.900430680208489, 0.496779210170732, 7.524616544037286)
    7  = @(1.459612070389937, 1.667794583268128, 3.900430680208489, 0.496779210170732, 7.524616544037286)
    8  = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 4.429675500412797)
    9  = @(0.3048080303191223, 0.3127903958511391, 0.8054896365839992, 0.496779210170732, 1.919867272924993)
    10 = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 4.429675500412797)
    11 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    12 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671)
}

foreach ($row in $values.Keys) {
    $v = $values[$row]
    $ws.Cells.Item($row, 2).Value = $v[0]
    $ws.Cells.Item($row, 3).Value = $v[1]
    $ws.Cells.Item($row, 4).Value = $v[2]
    $ws.Cells.Item($row, 5).Value = $v[3]
    $ws.Cells.Item($row, 7).Value = $v[4]
}
